# Select the "rank" worksheet (sheet2) and add the white-noise / script
# generated correlation column (C) next to the existing rank-correlation
# column (B).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rank")

$ws.Range("C2").Value = -0.0078
$ws.Range("C3").Value = -0.0527
$ws.Range("C4").Value = 0.0696
$ws.Range("C5").Value = 0.0139
$ws.Range("C6").Value = 0.0381
$ws.Range("C7").Value = 0.0172
$ws.Range("C8").Value = 0.0127
$ws.Range("C9").Value = 0.0078
$ws.Range("C10").Value = -0.0302
$ws.Range("C11").Value = 0.0192

# Match the resulting selection / active cell from the diff.
$ws.Range("C11").Select()

# Touch the page setup so the sheet carries explicit print settings, as in
# the committed workbook.
$ws.PageSetup.Orientation = 1
